# Updated cryptos list snapshot (prices / 1h volume deltas refreshed).
# Rows 30/31 and 39/40/41 were re-ranked, so Coin name / Link cells for
# those rows are rewritten too, in addition to Price (D) and Volume (E).
#
# D-column prices are plain numeric-looking text (e.g. "301.52",
# "0.0700") in the source workbook, so each Price write is wrapped with
# a temporary Text number format to stop Excel's COM layer from
# auto-coercing the string into a floating-point number (which would
# both change the stored value, e.g. 301.51999999999998, and drop
# trailing/leading zeros). The format is reset back to the sheet's
# default ("Normal" style) immediately after the write so no cell
# keeps a lingering custom number format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.186.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.91%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.311.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.57%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.83"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.505"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.14%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.508"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.80%  "
$ws.Range("E11").Value = "  +1.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.14"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.57%  "
$ws.Range("E13").Value = "  +3.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.75"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +15.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.80"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.675.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.71%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.377.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.811"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.113.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0907"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.35%  "
$ws.Range("E22").Value = "  +2.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "237.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.96%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.88%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.71%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.12"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +8.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0700"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.88%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.31%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.22%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.101"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.68%  "
$ws.Range("E42").Value = "  +0.66%  "
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.998.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0286"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "18.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.79%  "
$ws.Range("E48").Value = "  +4.73%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.97"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.539.45"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.32%  "
$ws.Range("E51").Value = "  +2.61%  "
